$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 75
$ws1.Range("F4").Value = 210
$ws1.Range("F5").Value = 1001
$ws1.Range("F6").Value = 206
$ws1.Range("F7").Value = 561
$ws1.Range("F8").Value = 108
$ws1.Range("F9").Value = 562
$ws1.Range("F10").Value = 563
$ws1.Range("F11").Value = 84
$ws1.Range("F12").Value = 38
$ws1.Range("F13").Value = 144
$ws1.Range("F14").Value = 220

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 6
$ws2.Range("F8").Value = 113

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6269
$ws3.Range("F3").Value = 769
$ws3.Range("F4").Value = 1885

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6269
$ws4.Range("F3").Value = 769
$ws4.Range("F4").Value = 1885
$ws4.Range("F7").Value = 75
$ws4.Range("F9").Value = 6
$ws4.Range("F11").Value = 210
$ws4.Range("F14").Value = 1001
$ws4.Range("F15").Value = 113
$ws4.Range("F16").Value = 206
$ws4.Range("F18").Value = 561
$ws4.Range("F20").Value = 108
$ws4.Range("F21").Value = 562
$ws4.Range("F23").Value = 563
$ws4.Range("F24").Value = 84
$ws4.Range("F27").Value = 38
$ws4.Range("F28").Value = 144
$ws4.Range("F34").Value = 220

$wb.Save()
